# fdr_worksheet works. Fields ready for validating
#
# Duplicate the existing "import_worksheet" sheet into a new sheet named
# "Procedure Based Requirements", placed right after the original, make it
# the active sheet/tab, and repair the sheet-scoped defined names
# (_FilterDatabase / Print_Titles / rngRequirements) so they point at the
# new sheet instead of (or in addition to) the original.

$wb = $excel.ActiveWorkbook

$orig = $wb.Worksheets.Item("import_worksheet")

# Copy the sheet, inserting the copy immediately after the original. This
# brings along cell values/formulas, styles, column widths, autofilter,
# page setup, header/footer, etc. "for free".
$orig.Copy($null, $orig)

# The copy lands right after $orig in tab order.
$newSheet = $wb.Worksheets.Item($orig.Index + 1)
$newSheet.Name = "Procedure Based Requirements"

# The copy operation clones the workbook-level "Print_Titles" defined name
# for the new sheet, but leaves its RefersTo pointing at the original
# sheet - repoint it at the new sheet.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Procedure Based Requirements!Print_Titles") {
        $n.RefersTo = "='Procedure Based Requirements'!`$1:`$1"
    }
}

# _FilterDatabase and rngRequirements aren't auto-created for the new
# sheet, so add sheet-scoped versions explicitly.
$newSheet.Names.Add("_xlnm._FilterDatabase", "='Procedure Based Requirements'!`$A`$1:`$O`$1")
$newSheet.Names.Add("rngRequirements", "=#REF!")

# Restore/adjust the remembered selections: the original sheet's cursor
# moves to I23, while the new sheet keeps the J7 selection it was copied
# with. Activating the new sheet last leaves it as the active tab.
$orig.Activate()
$orig.Range("I23").Select()

$newSheet.Activate()
$newSheet.Range("J7").Select()
